$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.258.89"
$ws.Range("E2").Value = "  -3.47%  "

$ws.Range("D3").Value = "'3.509.92"
$ws.Range("E3").Value = "  -4.77%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'581.99"
$ws.Range("E5").Value = "  -1.27%  "

$ws.Range("D6").Value = "'174.14"
$ws.Range("E6").Value = "  -3.91%  "

$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  +0.79%  "

$ws.Range("D8").Value = "'3.502.19"
$ws.Range("E8").Value = "  -4.79%  "

$ws.Range("D10").Value = "'0.190"
$ws.Range("E10").Value = "  -5.89%  "

$ws.Range("E11").Value = "  +5.78%  "

$ws.Range("E12").Value = "  -2.86%  "

$ws.Range("D13").Value = "'47.07"
$ws.Range("E13").Value = "  -6.03%  "

$ws.Range("E14").Value = "  -3.84%  "

$ws.Range("D15").Value = "'674.86"
$ws.Range("E15").Value = "  -1.70%  "

$ws.Range("D16").Value = "'4.070.31"
$ws.Range("E16").Value = "  -4.90%  "

$ws.Range("D17").Value = "'8.72"
$ws.Range("E17").Value = "  -3.64%  "

$ws.Range("D18").Value = "'69.209.95"
$ws.Range("E18").Value = "  -3.62%  "

$ws.Range("D19").Value = "'3.505.55"
$ws.Range("E19").Value = "  -4.81%  "

$ws.Range("E20").Value = "  -1.38%  "

$ws.Range("D21").Value = "'17.45"
$ws.Range("E21").Value = "  -3.94%  "

$ws.Range("D22").Value = "'11.19"
$ws.Range("E22").Value = "  -4.27%  "

$ws.Range("E23").Value = "  -4.37%  "

$ws.Range("D24").Value = "'16.17"
$ws.Range("E24").Value = "  -9.52%  "

$ws.Range("D25").Value = "'97.91"
$ws.Range("E25").Value = "  -6.12%  "

$ws.Range("E26").Value = "  -4.28%  "

$ws.Range("E27").Value = "  -0.62%  "

$ws.Range("E28").Value = "  +0.07%  "

$ws.Range("E29").Value = "  -6.78%  "

$ws.Range("D30").Value = "'9.46"
$ws.Range("E30").Value = "  -7.32%  "

$ws.Range("D31").Value = "'32.91"
$ws.Range("E31").Value = "  -7.14%  "

$ws.Range("E32").Value = "  -6.05%  "

$ws.Range("E33").Value = "  -7.89%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'7.33"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Value = "'1.36"
$ws.Range("E35").Value = "  -5.85%  "

$ws.Range("D36").Value = "'596.50"
$ws.Range("E36").Value = "  +5.14%  "

$ws.Range("E37").Value = "  -15.71%  "

$ws.Range("E38").Value = "  -3.72%  "

$ws.Range("E39").Value = "  -4.74%  "

$ws.Range("D40").Value = "'57.31"
$ws.Range("E40").Value = "  -3.81%  "

$ws.Range("E41").Value = "  +0.17%  "

$ws.Range("D42").Value = "'0.0439"
$ws.Range("E42").Value = "  -5.98%  "

$ws.Range("D43").Value = "'0.336"
$ws.Range("E43").Value = "  -4.95%  "

$ws.Range("E44").Value = "  -6.59%  "

$ws.Range("D45").Value = "'3.414.18"
$ws.Range("E45").Value = "  -9.74%  "

$ws.Range("D46").Value = "'33.44"
$ws.Range("E46").Value = "  -6.24%  "

$ws.Range("D47").Value = "'0.0₃0710"
$ws.Range("E47").Value = "  -9.10%  "

$ws.Range("D48").Value = "'2.90"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("D49").Value = "'2.61"
$ws.Range("E49").Value = "  -7.35%  "

$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("D51").Value = "'5.85"
$ws.Range("E51").Value = "  +19.42%  "
